$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 157, shifting existing rows 157:265 down to 158:266
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new record
$ws.Range("A157").Value = 11
$ws.Range("B157").Value = "Vega Monumental Concepción"
$ws.Range("C157").Value = "Bíobío"
$ws.Range("D157").Value = 44582
$ws.Range("E157").Value = 8
$ws.Range("F157").Value = "Fruta"
$ws.Range("G157").Value = 100101
$ws.Range("H157").Value = "Berries"
$ws.Range("I157").Value = 100112025
$ws.Range("J157").Value = "Frutilla"
$ws.Range("K157").Value = "Sin especificar"
$ws.Range("L157").Value = "Primera"
$ws.Range("M157").Value = 280
$ws.Range("N157").Value = 7000
$ws.Range("O157").Value = 7500
$ws.Range("P157").Value = 7268
$ws.Range("Q157").Value = "`$/caja 7 kilos"
$ws.Range("R157").Value = "Región del Maule"
$ws.Range("S157").Value = 1038
$ws.Range("T157").Value = 7
